# Auto-generated edit script: updates FFXIV leve-profit numbers on the Tiamat Profits workbook.
# Source: scheduled runner refresh of currentAveragePrice / Leve price / profit columns
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# --- ALC (sheet 1), row 11 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("I11").Value = 120
$ws.Range("H11").Value = 120
$ws.Range("K11").Value = 120
$ws.Range("M11").Value = 20

# --- ALC (sheet 1), row 18 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("M18").Value = 108.36667
$ws.Range("I18").Value = 175.63333
$ws.Range("H18").Value = 239.65625
$ws.Range("K18").Value = 175.63333

# --- ALC (sheet 1), row 33 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 526.3214
$ws.Range("K33").Value = 168.33333
$ws.Range("L33").Value = 2674.25
$ws.Range("N33").Value = -3132.25
$ws.Range("M33").Value = 60.66667000000001
$ws.Range("J33").Value = 2674.25
$ws.Range("I33").Value = 168.33333

# --- ALC (sheet 1), row 135 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("L135").Value = 500011470
$ws.Range("N135").Value = -500016540
$ws.Range("M135").Value = -1933.5
$ws.Range("J135").Value = 55556830
$ws.Range("I135").Value = 496.5
$ws.Range("H135").Value = 11111764
$ws.Range("K135").Value = 4468.5

# --- ALC (sheet 1), row 138 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("J138").Value = 2774.8142
$ws.Range("M138").Value = 1275.7999
$ws.Range("I138").Value = 1288.0667
$ws.Range("H138").Value = 2328.79
$ws.Range("K138").Value = 3864.2001
$ws.Range("L138").Value = 8324.442599999998
$ws.Range("N138").Value = -18604.4426

# --- ARM (sheet 2), row 32 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("I32").Value = 5792.298
$ws.Range("H32").Value = 6044.76
$ws.Range("K32").Value = 5792.298
$ws.Range("J32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("N32").Value = -10574
$ws.Range("M32").Value = -5505.298

# --- ARM (sheet 2), row 110 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("J110").Value = 1444
$ws.Range("M110").Value = 1518.25
$ws.Range("I110").Value = 526.75
$ws.Range("K110").Value = 526.75
$ws.Range("H110").Value = 796.5294
$ws.Range("L110").Value = 1444
$ws.Range("N110").Value = -5534

# --- BSM (sheet 3), row 132 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("N132").Value = -333383460
$ws.Range("J132").Value = 333373340
$ws.Range("H132").Value = 333373340
$ws.Range("L132").Value = 333373340

# --- BSM (sheet 3), row 141 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H141").Value = 49014.363
$ws.Range("L141").Value = 49014.363
$ws.Range("N141").Value = -59374.363
$ws.Range("J141").Value = 49014.363

# --- CRP (sheet 4), row 31 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("N31").Value = -16979
$ws.Range("M31").Value = -39736.73
$ws.Range("J31").Value = 16389
$ws.Range("I31").Value = 40031.73
$ws.Range("K31").Value = 40031.73
$ws.Range("H31").Value = 30359.705
$ws.Range("L31").Value = 16389

# --- CRP (sheet 4), row 34 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("N34").Value = -16793
$ws.Range("J34").Value = 16389
$ws.Range("M34").Value = -39829.73
$ws.Range("I34").Value = 40031.73
$ws.Range("H34").Value = 30359.705
$ws.Range("K34").Value = 40031.73
$ws.Range("L34").Value = 16389

# --- CRP (sheet 4), row 122 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("I122").Value = 1037.3334
$ws.Range("K122").Value = 3112.0002
$ws.Range("H122").Value = 1128.7142
$ws.Range("L122").Value = 3591.75
$ws.Range("N122").Value = -8491.75
$ws.Range("J122").Value = 1197.25
$ws.Range("M122").Value = -662.0001999999999

# --- CRP (sheet 4), row 132 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("N132").Value = -11079.8
$ws.Range("J132").Value = 2006.6
$ws.Range("M132").Value = -216683
$ws.Range("I132").Value = 73071
$ws.Range("H132").Value = 43460.832
$ws.Range("K132").Value = 219213
$ws.Range("L132").Value = 6019.799999999999

# --- CRP (sheet 4), row 134 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("J134").Value = 1900
$ws.Range("M134").Value = -1273.0587
$ws.Range("I134").Value = 1269.3529
$ws.Range("K134").Value = 3808.0587
$ws.Range("H134").Value = 1335.7368
$ws.Range("L134").Value = 5700
$ws.Range("N134").Value = -10770

# --- CUL (sheet 5), row 37 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H37").Value = 46800
$ws.Range("L37").Value = 140400
$ws.Range("N37").Value = -140624
$ws.Range("J37").Value = 46800

# --- CUL (sheet 5), row 131 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 135915.1
$ws.Range("L131").Value = 568811.34
$ws.Range("N131").Value = -578891.34
$ws.Range("J131").Value = 189603.78

# --- GSM (sheet 6), row 102 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("K102").Value = 10145.6
$ws.Range("L102").Value = 76082
$ws.Range("M102").Value = -8523.6
$ws.Range("N102").Value = -79326
$ws.Range("J102").Value = 76082
$ws.Range("I102").Value = 10145.6
$ws.Range("H102").Value = 24026.947

# --- GSM (sheet 6), row 122 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("I122").Value = 2632
$ws.Range("K122").Value = 7896
$ws.Range("H122").Value = 3158.111
$ws.Range("L122").Value = 14998.5
$ws.Range("N122").Value = -19898.5
$ws.Range("J122").Value = 4999.5
$ws.Range("M122").Value = -5446

# --- GSM (sheet 6), row 132 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("N132").Value = -1012287.02
$ws.Range("J132").Value = 335742.34
$ws.Range("M132").Value = -4483.499899999999
$ws.Range("I132").Value = 2337.8333
$ws.Range("H132").Value = 113472.664
$ws.Range("K132").Value = 7013.499899999999
$ws.Range("L132").Value = 1007227.02

# --- LTW (sheet 7), row 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("K7").Value = 2428.5652
$ws.Range("L7").Value = 13668.889
$ws.Range("N7").Value = -13892.889
$ws.Range("M7").Value = -2316.5652
$ws.Range("J7").Value = 13668.889
$ws.Range("I7").Value = 2428.5652
$ws.Range("H7").Value = 5589.9062

# --- LTW (sheet 7), row 126 ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("I126").Value = 2428.5652
$ws.Range("H126").Value = 5589.9062
$ws.Range("K126").Value = 7285.6956
$ws.Range("L126").Value = 41006.667
$ws.Range("J126").Value = 13668.889
$ws.Range("N126").Value = -45946.667
$ws.Range("M126").Value = -4815.6956

# --- LTW (sheet 7), row 132 ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("N132").Value = -1381123.1
$ws.Range("J132").Value = 458687.7
$ws.Range("M132").Value = -179476.298
$ws.Range("I132").Value = 60668.766
$ws.Range("H132").Value = 217033.34
$ws.Range("K132").Value = 182006.298
$ws.Range("L132").Value = 1376063.1

# --- LTW (sheet 7), row 136 ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("K136").Value = 814460.7000000001
$ws.Range("L136").Value = 6423.75
$ws.Range("N136").Value = -11523.75
$ws.Range("J136").Value = 2141.25
$ws.Range("M136").Value = -811910.7000000001
$ws.Range("I136").Value = 271486.9
$ws.Range("H136").Value = 223603.22

# --- WVR (sheet 8), row 96 ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("N96").Value = -6679.3333
$ws.Range("M96").Value = -878.5
$ws.Range("J96").Value = 3933.3333
$ws.Range("I96").Value = 2251.5
$ws.Range("H96").Value = 3260.6
$ws.Range("K96").Value = 2251.5
$ws.Range("L96").Value = 3933.3333

# --- WVR (sheet 8), row 107 ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("I107").Value = 323
$ws.Range("H107").Value = 323
$ws.Range("K107").Value = 969
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 951
$ws.Range("N107").ClearContents()

# --- WVR (sheet 8), row 113 ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("I113").Value = 362.6154
$ws.Range("H113").Value = 638.0357
$ws.Range("K113").Value = 1087.8462
$ws.Range("L113").Value = 2630.20002
$ws.Range("J113").Value = 876.73334
$ws.Range("N113").Value = -6970.20002
$ws.Range("M113").Value = 1082.1538

# --- WVR (sheet 8), row 136 ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("K136").Value = 6127828.199999999
$ws.Range("L136").Value = 2734993.68
$ws.Range("N136").Value = -2740093.68
$ws.Range("J136").Value = 911664.5600000001
$ws.Range("M136").Value = -6125278.199999999
$ws.Range("I136").Value = 2042609.4
$ws.Range("H136").Value = 1772166
